# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.856.01"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.556.79"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'315.45"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'96.36"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'35.45"
$ws.Range("E10").Value = "  -3.45%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'7.42"
$ws.Range("E12").Value = "  -3.10%  "
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("D14").Value = "2.949.63"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "2.568.83"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "'15.07"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").Value = "'0.840"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "42.919.63"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'6.81"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "'12.50"
$ws.Range("E20").Value = "  -4.67%  "
$ws.Range("D21").Value = "0.0₃0957"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'69.17"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").Value = "'251.74"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'2.94"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").Value = "'26.76"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'2.43"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "'39.81"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'10.15"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "'5.78"
$ws.Range("E31").Value = "  -5.67%  "
$ws.Range("D32").Value = "'154.23"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "'3.40"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.70"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0803"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "'2.11"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "'19.10"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'2.44"
$ws.Range("E39").Value = "  +6.64%  "
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").Value = "'22.52"
$ws.Range("E41").Value = "  -7.80%  "
$ws.Range("D42").Value = "'3.92"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'3.24"
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("D46").Value = "2.000.16"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "'9.00"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").Value = "'83.05"
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("D49").Value = "2.804.39"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'73.89"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'103.44"
$ws.Range("E51").Value = "  +0.26%  "
